$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells that would otherwise be parsed as numbers
$textForceCells = @("D4", "D5", "D6", "D8", "D9", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D21", "D22", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '29.189.05'
$ws.Range("E2").Value = '  +0.00%  '

$ws.Range("D3").Value = '1.839.27'

$ws.Range("D4").Value = '1.000'

$ws.Range("D5").Value = '242.60'
$ws.Range("E5").Value = '  +0.76%  '

$ws.Range("D6").Value = '0.6621'
$ws.Range("E6").Value = '  -1.43%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").Value = '44.72'
$ws.Range("E8").Value = '  +6.39%  '

$ws.Range("D9").Value = '0.07435'
$ws.Range("E9").Value = '  +0.23%  '

$ws.Range("E10").Value = '  -0.01%  '

$ws.Range("D11").Value = '23.21'
$ws.Range("E11").Value = '  +1.64%  '

$ws.Range("D12").Value = '0.07732'
$ws.Range("E12").Value = '  +0.14%  '

$ws.Range("D13").Value = '1.840.49'
$ws.Range("E13").Value = '  -0.27%  '

$ws.Range("D14").Value = '5.011'
$ws.Range("E14").Value = '  -0.01%  '

$ws.Range("D15").Value = '0.6727'
$ws.Range("E15").Value = '  -0.62%  '

$ws.Range("D16").Value = '83.08'
$ws.Range("E16").Value = '  -3.53%  '

$ws.Range("D17").Value = '6.166'
$ws.Range("E17").Value = '  +0.44%  '

$ws.Range("D18").Value = '0.000008713'
$ws.Range("E18").Value = '  +4.92%  '

$ws.Range("D19").Value = '29.200.66'
$ws.Range("E19").Value = '  +0.02%  '

$ws.Range("D20").Value = '2.091.73'
$ws.Range("E20").Value = '  +0.42%  '

$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").Value = '12.52'
$ws.Range("E21").Value = '  -0.17%  '

$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").Value = '226.10'
$ws.Range("E22").Value = '  -1.09%  '

$ws.Range("E23").Value = '  -0.01%  '

$ws.Range("D24").Value = '7.169'
$ws.Range("E24").Value = '  -0.40%  '

$ws.Range("E25").Value = '  +0.02%  '

$ws.Range("D26").Value = '158.80'
$ws.Range("E26").Value = '  -1.30%  '

$ws.Range("D27").Value = '8.608'
$ws.Range("E27").Value = '  -0.91%  '

$ws.Range("D28").Value = '0.1392'
$ws.Range("E28").Value = '  -1.05%  '

$ws.Range("D29").Value = '18.03'

$ws.Range("D30").Value = '1.511'
$ws.Range("E30").Value = '  +0.13%  '

$ws.Range("D31").Value = '4.137'
$ws.Range("E31").Value = '  -0.91%  '

$ws.Range("D32").Value = '1.209'
$ws.Range("E32").Value = '  +1.11%  '

$ws.Range("D33").Value = '4.041'
$ws.Range("E33").Value = '  -0.81%  '

$ws.Range("E34").Value = '  +0.97%  '

$ws.Range("D35").Value = '1.854'
$ws.Range("E35").Value = '  -1.84%  '

$ws.Range("D36").Value = '0.7484'
$ws.Range("E36").Value = '  -1.34%  '

$ws.Range("D37").Value = '1.160'
$ws.Range("E37").Value = '  +2.01%  '

$ws.Range("D38").Value = '2.653'
$ws.Range("E38").Value = '  -1.29%  '

$ws.Range("D39").Value = '1.300.42'
$ws.Range("E39").Value = '  -2.18%  '

$ws.Range("D40").Value = '0.01795'
$ws.Range("E40").Value = '  -0.41%  '

$ws.Range("D41").Value = '2.761'
$ws.Range("E41").Value = '  +1.04%  '

$ws.Range("D42").Value = '6.353'
$ws.Range("E42").Value = '  +6.47%  '

$ws.Range("D43").Value = '0.9033'

$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  -0.12%  '

$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '103.82'
$ws.Range("E45").Value = '  +0.37%  '

$ws.Range("B46").Value = 'XinFinNetwork'
$ws.Range("C46").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D46").Value = '0.08242'
$ws.Range("E46").Value = '  +3.13%  '

$ws.Range("D47").Value = '1.988.97'
$ws.Range("E47").Value = '  +1.13%  '

$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '65.10'
$ws.Range("E48").Value = '  +1.55%  '

$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.00000000122'
$ws.Range("E49").Value = '  -2.01%  '

$ws.Range("D50").Value = '0.5142'
$ws.Range("E50").Value = '  -0.51%  '

$ws.Range("D51").Value = '1.752'
$ws.Range("E51").Value = '  -1.10%  '

# Restore default style on cells where we forced text number format,
# so the cell style matches the original (no explicit style index)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}